# Generate Report for Handoff
#
# The handoff/report generation run completed for the four files that were
# still "Ready for handoff" (06134032-..., 3c1e5da5-..., 6c37ebc1-...,
# 81ce6a9e-...). For each language sheet this:
#   - bumps the Priority for those rows from "low" to "ht"
#   - records the new "Latest Handoff Datetime" for those rows
# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff timestamp for the same rows, so it is updated to match.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-24 00:30:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-24 00:30:36"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-24 00:30:36"
